$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.268.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.420.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.91%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.417.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("E11").Value = "  -2.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.859.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.098.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.419.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("E25").Value = "  -3.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "578.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0955"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.540.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  +3.29%  "
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "152.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("E42").Value = "  -4.57%  "
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "149.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("E46").Value = "  +1.44%  "
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.596"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0921"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.78%  "
$ws.Range("E51").Value = "  +2.08%  "
